# Auto-generated Excel COM-interop script
# Applies a scheduled market-price data refresh to the Lamia_Profits workbook:
# updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) across
# several rows on the ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 588.41174
$ws.Range("I11").Value = 588.41174
$ws.Range("K11").Value = 588.41174
$ws.Range("M11").Value = -448.41174
$ws.Range("H18").Value = 1942.5
$ws.Range("I18").Value = 923.3333
$ws.Range("K18").Value = 923.3333
$ws.Range("M18").Value = -639.3333
$ws.Range("H28").Value = 1764.9286
$ws.Range("I28").Value = 391.0476
$ws.Range("J28").Value = 5886.5713
$ws.Range("K28").Value = 391.0476
$ws.Range("L28").Value = 5886.5713
$ws.Range("M28").Value = 93.95240000000001
$ws.Range("N28").Value = -6856.5713
$ws.Range("H64").Value = 19003
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 19003
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H99").Value = 6842.6665
$ws.Range("I99").Value = 6842.6665
$ws.Range("K99").Value = 20527.9995
$ws.Range("M99").Value = -19029.9995
$ws.Range("H107").Value = 942.1875
$ws.Range("I107").Value = 690.5
$ws.Range("K107").Value = 690.5
$ws.Range("M107").Value = 1229.5
$ws.Range("H138").Value = 7725.9614
$ws.Range("J138").Value = 7627.12
$ws.Range("L138").Value = 22881.36
$ws.Range("N138").Value = -33161.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 382.66666
$ws.Range("J5").Value = 499
$ws.Range("L5").Value = 499
$ws.Range("N5").Value = -723
$ws.Range("H32").Value = 14306.19
$ws.Range("I32").Value = 11021.5
$ws.Range("J32").Value = 80000
$ws.Range("K32").Value = 11021.5
$ws.Range("L32").Value = 80000
$ws.Range("M32").Value = -10734.5
$ws.Range("N32").Value = -80574
$ws.Range("H44").Value = 14000
$ws.Range("I44").Value = 14000
$ws.Range("K44").Value = 14000
$ws.Range("M44").Value = -13512
$ws.Range("H55").Value = 10341.333
$ws.Range("I55").Value = 10341.333
$ws.Range("K55").Value = 10341.333
$ws.Range("M55").Value = -10026.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 382.66666
$ws.Range("J4").Value = 499
$ws.Range("L4").Value = 499
$ws.Range("N4").Value = -729
$ws.Range("H22").Value = 1975.7778
$ws.Range("I22").Value = 407.4
$ws.Range("K22").Value = 407.4
$ws.Range("M22").Value = -234.4
$ws.Range("H80").Value = 1574.1111
$ws.Range("I80").Value = 300
$ws.Range("J80").Value = 1938.1428
$ws.Range("K80").Value = 300
$ws.Range("L80").Value = 1938.1428
$ws.Range("M80").Value = 698
$ws.Range("N80").Value = -3934.1428
$ws.Range("H83").Value = 1574.1111
$ws.Range("I83").Value = 300
$ws.Range("J83").Value = 1938.1428
$ws.Range("K83").Value = 1500
$ws.Range("L83").Value = 9690.714
$ws.Range("M83").Value = 3492
$ws.Range("N83").Value = -19674.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 3050
$ws.Range("I25").Value = 1100
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 1100
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = -926
$ws.Range("N25").Value = -5348
$ws.Range("H141").Value = 691174.75
$ws.Range("J141").Value = 785628.3
$ws.Range("L141").Value = 785628.3
$ws.Range("N141").Value = -795988.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 460270.38
$ws.Range("I80").Value = 1003734.4
$ws.Range("J80").Value = 7383.6665
$ws.Range("K80").Value = 1003734.4
$ws.Range("L80").Value = 7383.6665
$ws.Range("M80").Value = -1002736.4
$ws.Range("N80").Value = -9379.666499999999
$ws.Range("H83").Value = 460270.38
$ws.Range("I83").Value = 1003734.4
$ws.Range("J83").Value = 7383.6665
$ws.Range("K83").Value = 5018672
$ws.Range("L83").Value = 36918.3325
$ws.Range("M83").Value = -5013680
$ws.Range("N83").Value = -46902.3325
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H134").Value = 63526.855
$ws.Range("J134").Value = 63526.855
$ws.Range("L134").Value = 190580.565
$ws.Range("N134").Value = -195650.565
$ws.Range("H136").Value = 47140.855
$ws.Range("J136").Value = 47140.855
$ws.Range("L136").Value = 141422.565
$ws.Range("N136").Value = -146522.565

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 14650.75
$ws.Range("I22").Value = 1766.5
$ws.Range("J22").Value = 22381.3
$ws.Range("K22").Value = 1766.5
$ws.Range("L22").Value = 22381.3
$ws.Range("M22").Value = -1471.5
$ws.Range("N22").Value = -22971.3
$ws.Range("H23").Value = 13927.5
$ws.Range("I23").Value = 13927.5
$ws.Range("K23").Value = 13927.5
$ws.Range("M23").Value = -13697.5
$ws.Range("H27").Value = 14650.75
$ws.Range("I27").Value = 1766.5
$ws.Range("J27").Value = 22381.3
$ws.Range("K27").Value = 1766.5
$ws.Range("L27").Value = 22381.3
$ws.Range("M27").Value = -1659.5
$ws.Range("N27").Value = -22595.3
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H42").Value = 33341.668
$ws.Range("I42").Value = 20025
$ws.Range("K42").Value = 20025
$ws.Range("M42").Value = -19462
$ws.Range("H49").Value = 33341.668
$ws.Range("I49").Value = 20025
$ws.Range("K49").Value = 20025
$ws.Range("M49").Value = -19878
$ws.Range("H53").Value = 19797.334
$ws.Range("I53").Value = 9696
$ws.Range("J53").Value = 40000
$ws.Range("K53").Value = 9696
$ws.Range("L53").Value = 40000
$ws.Range("M53").Value = -9178
$ws.Range("N53").Value = -41036
$ws.Range("H136").Value = 4955.8237
$ws.Range("I136").Value = 3737.5833
$ws.Range("J136").Value = 7879.6
$ws.Range("K136").Value = 11212.7499
$ws.Range("L136").Value = 23638.8
$ws.Range("M136").Value = -8662.749899999999
$ws.Range("N136").Value = -28738.8
$ws.Range("H137").Value = 69002.14
$ws.Range("J137").Value = 69002.14
$ws.Range("L137").Value = 69002.14
$ws.Range("N137").Value = -79202.14

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1809.409
$ws.Range("I107").Value = 1782.5333
$ws.Range("J107").Value = 1867
$ws.Range("K107").Value = 5347.5999
$ws.Range("L107").Value = 5601
$ws.Range("M107").Value = -3427.5999
$ws.Range("N107").Value = -9441
